# Project Sample Project is saved. The "Rules" sheet's last rule row (R40)
# has its "To" label (B11) changed from "R40" to "1".
# The leading apostrophe forces Excel to store the numeric-looking literal
# as text (matching the existing column's text values R10/R20/R30/...)
# instead of silently converting it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("B11").Value = "'1"
